# Estadisticos Segundo Parcial 23 Mayo
# Update the partial/final statistics sheets and clear out the
# "Rescatables" (resit-eligible students) listing.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: "Estadisticos 1P"
# Row 2 (Formación socioemocional II / 2ASV) statistics reset to 0
# and the computed Por_Apro / Promedio cells are cleared.
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Estadisticos 1P")
$ws1.Range("C2").Value = 0
$ws1.Range("E2").Value = 0
$ws1.Range("F2").Value = 0
$ws1.Range("G2").ClearContents()
$ws1.Range("H2").ClearContents()

# ---------------------------------------------------------------
# Sheet 2: "Estadisticos 2P"
# Row 2 (Formación socioemocional II) statistics reset to 0 and the
# Por_Apro cell is cleared.
# Row 3 (Lengua y comunicación II) now has real figures calculated.
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")
$ws2.Range("C2").Value = 0
$ws2.Range("D2").Value = 0
$ws2.Range("E2").Value = 0
$ws2.Range("G2").ClearContents()

$ws2.Range("D3").Value = 0
$ws2.Range("E3").Value = 2
$ws2.Range("F3").Value = 9
$ws2.Range("G3").Value = 81.81999999999999
$ws2.Range("H3").Value = 7.8

# ---------------------------------------------------------------
# Sheet 3: "Estadisticos Final"
# Row 2 (Formación socioemocional II) statistics reset to 0 and the
# Por_Apro cell is cleared.
# Row 3 (Lengua y comunicación II) average is updated.
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Estadisticos Final")
$ws3.Range("C2").Value = 0
$ws3.Range("D2").Value = 0
$ws3.Range("G2").ClearContents()

$ws3.Range("H3").Value = 8.4

# ---------------------------------------------------------------
# Sheet 4: "Rescatables"
# Remove all the listed resit-eligible students, keeping only the
# header row.
# ---------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Rescatables")
$ws4.Rows("2:9").Delete()
